# Change os to windows
# Rows 2 and 3 (column A, "Categories") currently hold "OS"; the commit
# changes them to the new value "Windows" (a new shared string). The
# active selection also moves from C11 to A3, matching where the edit
# was made in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Windows"
$ws.Range("A3").Value = "Windows"

$ws.Range("A3").Select()
